$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.628.85'
$ws.Range("E2").Value = '  -1.64%  '

$ws.Range("D3").Value = '1.590.49'
$ws.Range("E3").Value = '  -2.14%  '

$ws.Range("E4").Value = '  -0.32%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.75%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.512'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.12%  '

$ws.Range("E7").Value = '  -0.17%  '

$ws.Range("E8").Value = '  -2.30%  '

$ws.Range("E9").Value = '  -1.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.67'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0834'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.56%  '

$ws.Range("D12").Value = '1.813.46'
$ws.Range("E12").Value = '  -2.54%  '

$ws.Range("D13").Value = '1.596.14'
$ws.Range("E13").Value = '  -2.30%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.24%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.526'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.59%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.75'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.66%  '

$ws.Range("D17").Value = '26.628.33'
$ws.Range("E17").Value = '  -1.83%  '

$ws.Range("D18").Value = '0.0₃0727'
$ws.Range("E18").Value = '  -1.53%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '208.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.53%  '

$ws.Range("E20").Value = '  -0.02%  '

$ws.Range("E21").Value = '  -1.90%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.25'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.77%  '

$ws.Range("E23").Value = '  -2.15%  '

$ws.Range("E24").Value = '  -0.99%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.78'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.56%  '

$ws.Range("E26").Value = '  -0.24%  '

$ws.Range("E27").Value = '  -0.09%  '

$ws.Range("E28").Value = '  -3.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.31'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.71%  '

$ws.Range("E30").Value = '  +0.44%  '

$ws.Range("E31").Value = '  -1.81%  '

$ws.Range("E32").Value = '  -3.57%  '

$ws.Range("E33").Value = '  +22.44%  '

$ws.Range("E34").Value = '  -2.33%  '

$ws.Range("D35").Value = '1.318.14'
$ws.Range("E35").Value = '  -1.21%  '

$ws.Range("E36").Value = '  -4.33%  '

$ws.Range("E37").Value = '  -2.10%  '

$ws.Range("E38").Value = '  -1.86%  '

$ws.Range("E39").Value = '  -2.02%  '

$ws.Range("E40").Value = '  -0.06%  '

$ws.Range("E41").Value = '  +3.16%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.790'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.74%  '

$ws.Range("E43").Value = '  -3.58%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.30'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.42%  '

$ws.Range("D45").Value = '1.726.64'
$ws.Range("E45").Value = '  -2.41%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.99'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.52%  '

$ws.Range("E47").Value = '  +0.14%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.840'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.21%  '

$ws.Range("E49").Value = '  -0.24%  '

$ws.Range("E50").Value = '  -0.78%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.52'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.39%  '
